# Rename rouge1 -> rougel headers (mean/std) and update associated metric
# values across all four worksheets, per the "feat: change to correct
# rouge" commit.

$wb = $excel.ActiveWorkbook

# New header text + new numeric-looking text values (kept as Text so the
# comma-decimal strings are preserved verbatim, matching the workbook's
# existing string-cell storage).
$headerMap = @{
    "B1" = "mean_rougel_f"
    "C1" = "std_rougel_f"
}

# sheet name -> row -> @(B-value, C-value)
$valueMap = @{
    "wanetal_results_report" = @{
        2  = @("0,4135880628838103", "0,273585971981938")
        3  = @("0,40710156602483993", "0,26825571579833285")
        4  = @("0,36758372684168195", "0,2270243545339394")
        5  = @("0,4546510352479431", "0,30207758652510003")
        6  = @("0,41749451748297195", "0,2672183291578415")
        7  = @("0,4568390471703565", "0,30145102617980485")
        8  = @("0,36843757548969347", "0,22670339508873064")
        9  = @("0,4103334471824482", "0,27039355768925033")
        10 = @("0,36806787051591566", "0,23316500047691527")
        11 = @("0,456952155105607", "0,3023002624791239")
        12 = @("0,4161868475004489", "0,26752921789007095")
        13 = @("0,4099304562893027", "0,2696637576467272")
        14 = @("0,45669792766638645", "0,3022918320940732")
        15 = @("0,40384354391282895", "0,254215274223736")
        16 = @("0,3677914373223191", "0,22829889400795722")
    }
    "codexglue_results_report" = @{
        2 = @("0,3304149838939274", "0,19603554838058318")
        3 = @("0,3482884438287336", "0,1948751835361479")
        4 = @("0,32734240196930287", "0,18506659375230844")
        5 = @("0,3032659666170646", "0,15042429719487305")
        6 = @("0,35020986004032917", "0,1945378727069066")
    }
    "huetal_results_report" = @{
        2  = @("0,4428271805975517", "0,2766894351597228")
        3  = @("0,44077007868241846", "0,2618179410112169")
        4  = @("0,44110618878971314", "0,2636937177731049")
        5  = @("0,4933986961021123", "0,3012508235213127")
        6  = @("0,44266157804112644", "0,2759802534357222")
        7  = @("0,4903718104370656", "0,2897250214134543")
        8  = @("0,4921043414011757", "0,30010289686379676")
        9  = @("0,4897661797317905", "0,30313109480237155")
        10 = @("0,4113543177481951", "0,24593006779025553")
        11 = @("0,47719847420166717", "0,2764852119782333")
        12 = @("0,487165872049073", "0,2876688271937383")
        13 = @("0,4179078949438981", "0,2373207047777384")
        14 = @("0,4170620364778556", "0,2358367382523934")
        15 = @("0,41665098557853886", "0,23726232156539284")
        16 = @("0,49029780455272187", "0,3019908964766978")
    }
    "java_codexglue_results_report" = @{
        2 = @("0,3827609327428907", "0,22172419618991876")
        3 = @("0,3776550049208594", "0,22354942556658733")
        4 = @("0,3637469843043001", "0,20617041731096186")
        5 = @("0,40478078157851205", "0,2287514226159065")
        6 = @("0,4030341465129616", "0,23140391793432183")
    }
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $valueMap.ContainsKey($name)) {
        continue
    }

    # Update the headers.
    $ws.Range("B1").Value = $headerMap["B1"]
    $ws.Range("C1").Value = $headerMap["C1"]

    $rows = $valueMap[$name]
    foreach ($r in $rows.Keys) {
        $pair = $rows[$r]

        # Force Text format first so these comma-decimal numeric-looking
        # strings are not re-interpreted/mangled as numbers (Excel would
        # otherwise treat the comma as a thousands separator).
        $bCell = $ws.Range("B$r")
        $cCell = $ws.Range("C$r")
        $bCell.NumberFormat = "@"
        $cCell.NumberFormat = "@"
        $bCell.Value = $pair[0]
        $cCell.Value = $pair[1]
    }
}
